$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1) Locate the paragraph that ends with "...Como desde la experiencia
#    se construye una práctica y un proyecto. " (section 3.2) — this is
#    where the new italic run "CONCEPTO POLITICO DE TERNURA" is added.
# --------------------------------------------------------------------
$anchorText = "Como desde la experiencia se construye una práctica y un proyecto. "
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Contains($anchorText)) {
        $targetParaIndex = $i
        break
    }
}
if ($targetParaIndex -eq -1) {
    throw "Could not find target paragraph"
}

# --------------------------------------------------------------------
# 2) Remove the existing (hidden) "_GoBack" bookmark — it currently
#    sits at the very end of the document and needs to move.
# --------------------------------------------------------------------
try {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
} catch {
    # no pre-existing _GoBack bookmark — nothing to remove
}

# --------------------------------------------------------------------
# 3) Append "CONCEPTO POLITICO DE TERNURA" right after the existing
#    italic run, inheriting its formatting (italic + es-AR) by growing
#    the run via Find/Replace, then splitting it back into two runs
#    (the appended text becomes a run of its own, matching the target
#    markup) by nudging a character-formatting property.
# --------------------------------------------------------------------
$newWords = "CONCEPTO POLITICO DE TERNURA"
$oldText = $anchorText
$newText = $oldText + $newWords
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2) | Out-Null

$p = $d.Paragraphs($targetParaIndex)
$paraRange = $p.Range
$paraText = $paraRange.Text
$wordsIdx = $paraText.IndexOf($newWords)
$segStart = $paraRange.Start + $wordsIdx
$segEnd = $segStart + $newWords.Length
$newRunRange = $d.Range($segStart, $segEnd)
# Toggling Bold off again is a formatting no-op but forces the new text
# to be serialized as its own <w:r> (with its own, explicit <w:rPr>)
# instead of being merged back into the preceding run.
$newRunRange.Font.Bold = $true
$newRunRange.Font.Bold = $false

# --------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark immediately after the new run,
#    i.e. at the (now shifted) end of the paragraph. Adding a
#    zero-length bookmark directly at a paragraph-end offset isn't
#    reliable, so a placeholder character is inserted, bookmarked
#    together with it, and then removed again — which correctly
#    collapses the bookmark back down to zero width in place.
# --------------------------------------------------------------------
$p = $d.Paragraphs($targetParaIndex)
$endPos = $p.Range.End - 1
$p.Range.InsertAfter("X")
$placeholderRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange) | Out-Null
$placeholderRange2 = $d.Range($endPos, $endPos + 1)
$placeholderRange2.Delete()
